$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updated odds values
$ws.Range("G2").Value = 1.75
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 4.85
$ws.Range("J2").Value = 2.35
$ws.Range("K2").Value = 2.02
$ws.Range("L2").Value = 5.1
$ws.Range("M2").Value = 1.27
$ws.Range("N2").Value = 3.1
$ws.Range("O2").Value = 1.8
$ws.Range("P2").Value = 1.8
$ws.Range("Q2").Value = 2.82
$ws.Range("R2").Value = 1.32
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 1.7
$ws.Range("V2").Value = 1.93
$ws.Range("W2").Value = 7.1
$ws.Range("X2").Value = 8.5
$ws.Range("Y2").Value = 7.8
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 13.5
$ws.Range("AB2").Value = 23
$ws.Range("AC2").Value = 9.5
$ws.Range("AD2").Value = 6.3
$ws.Range("AE2").Value = 13.5
$ws.Range("AF2").Value = 60
$ws.Range("AG2").Value = 450
$ws.Range("AH2").Value = 13
$ws.Range("AI2").Value = 30
$ws.Range("AJ2").Value = 15
$ws.Range("AK2").Value = 100
$ws.Range("AL2").Value = 50

# Row 9 updated odds values
$ws.Range("G9").Value = 2.95
$ws.Range("H9").Value = 3.7
$ws.Range("I9").Value = 2.07
$ws.Range("J9").Value = 3.4
$ws.Range("K9").Value = 2.27
$ws.Range("L9").Value = 2.62
$ws.Range("M9").Value = 1.2
$ws.Range("N9").Value = 3.6
$ws.Range("O9").Value = 1.6
$ws.Range("P9").Value = 2.05
$ws.Range("Q9").Value = 2.42
$ws.Range("R9").Value = 1.44
$ws.Range("U9").Value = 1.55
$ws.Range("V9").Value = 2.15
$ws.Range("W9").Value = 11.75
$ws.Range("X9").Value = 17.5
$ws.Range("Y9").Value = 10.75
$ws.Range("Z9").Value = 37
$ws.Range("AA9").Value = 23
$ws.Range("AC9").Value = 13.5
$ws.Range("AD9").Value = 7.3
$ws.Range("AF9").Value = 45
$ws.Range("AG9").Value = 300
$ws.Range("AI9").Value = 11.5
$ws.Range("AJ9").Value = 8.75
$ws.Range("AK9").Value = 20
$ws.Range("AL9").Value = 15.5
$ws.Range("AM9").Value = 22
